$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9633786082267761
$ws.Range("B1").Value = 1.610763669013977
$ws.Range("D1").Value = 1.785297989845276
$ws.Range("E1").Value = 1.066914558410645
